$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as Text before assignment so that
# numeric-looking strings (e.g. "0.0000258", "70.950.86") are preserved
# exactly as text instead of being auto-converted to numbers.
$cellValues = @{
    "D2" = "70.950.86"
    "E2" = "  +2.53%  "
    "D3" = "3.808.92"
    "E4" = "  -0.11%  "
    "D5" = "698.83"
    "E5" = "  +10.35%  "
    "D6" = "172.67"
    "E6" = "  +3.59%  "
    "D7" = "3.808.33"
    "E7" = "  +0.94%  "
    "E8" = "  +0.00%  "
    "D9" = "0.525"
    "E9" = "  +0.72%  "
    "E10" = "  +2.61%  "
    "D11" = "7.52"
    "E11" = "  +11.56%  "
    "E12" = "  +0.53%  "
    "D13" = "0.0000258"
    "E13" = "  +5.76%  "
    "D14" = "36.33"
    "E14" = "  +3.61%  "
    "D15" = "4.451.16"
    "E15" = "  +0.93%  "
    "D16" = "3.816.02"
    "E16" = "  +1.05%  "
    "D17" = "71.063.60"
    "E17" = "  +2.66%  "
    "D18" = "17.81"
    "E18" = "  +1.08%  "
    "D19" = "7.22"
    "E19" = "  +2.76%  "
    "E20" = "  +0.14%  "
    "D21" = "11.31"
    "E21" = "  +18.47%  "
    "D22" = "480.36"
    "E22" = "  +3.55%  "
    "E23" = "  +1.47%  "
    "D24" = "83.94"
    "E24" = "  +1.48%  "
    "D25" = "0.0000145"
    "E25" = "  +0.83%  "
    "D26" = "12.39"
    "E26" = "  +2.52%  "
    "D28" = "10.45"
    "E28" = "  +3.60%  "
    "D29" = "3.959.87"
    "E29" = "  +0.93%  "
    "E30" = "  -0.09%  "
    "D31" = "3.09"
    "E31" = "  +14.77%  "
    "D32" = "2.32"
    "E32" = "  +0.11%  "
    "D33" = "7.55"
    "E33" = "  +6.31%  "
    "D34" = "29.58"
    "E34" = "  +3.81%  "
    "E35" = "  +5.87%  "
    "D36" = "9.23"
    "E36" = "  +2.74%  "
    "E37" = "  +0.02%  "
    "D38" = "3.758.29"
    "E38" = "  +0.85%  "
    "E39" = "  +1.35%  "
    "D40" = "3.49"
    "E40" = "  +6.18%  "
    "E41" = "  +3.11%  "
    "B42" = "Stacks"
    "C42" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D42" = "2.23"
    "E42" = "  +13.31%  "
    "B43" = "FLOKI"
    "C43" = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
    "D43" = "0.000332"
    "E43" = "  +24.09%  "
    "D44" = "0.971"
    "E44" = "  +1.09%  "
    "E45" = "  -0.03%  "
    "E46" = "  -0.01%  "
    "D47" = "45.42"
    "E47" = "  +5.13%  "
    "D48" = "160.90"
    "E48" = "  +1.73%  "
    "D49" = "49.15"
    "E49" = "  +5.08%  "
    "E50" = "  -1.02%  "
    "D51" = "0.300"
    "E51" = "  +1.10%  "
}

foreach ($addr in $cellValues.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellValues[$addr]
}